$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 29.53999999999999
$ws.Range("C2").Value = 29.42570304870605
$ws.Range("D2").Value = -0.1142969512939374
$ws.Range("E2").Value = 0.01306379307508869
$ws.Range("C3").Value = 29.50772666931152
$ws.Range("D3").Value = -0.04227333068847372
$ws.Range("E3").Value = 0.001787034487497054
$ws.Range("C4").Value = 29.72051811218262
$ws.Range("D4").Value = -0.02948188781738281
$ws.Range("E4").Value = 0.0008691817092767451
$ws.Range("C5").Value = 29.87379264831543
$ws.Range("D5").Value = 0.03379264831542628
$ws.Range("E5").Value = 0.001141943080170082
$ws.Range("C6").Value = 29.8944149017334
$ws.Range("D6").Value = 0.08441490173339616
$ws.Range("E6").Value = 0.007125875634658931
$ws.Range("C7").Value = 29.93996238708496
$ws.Range("D7").Value = 0.01996238708495923
$ws.Range("E7").Value = 0.0003984968981297471
$ws.Range("C8").Value = 30.07316398620605
$ws.Range("D8").Value = 0.09316398620605071
$ws.Range("E8").Value = 0.008679528325801206
$ws.Range("B9").Value = 30.03999999999999
$ws.Range("C9").Value = 30.21114921569824
$ws.Range("D9").Value = 0.1711492156982501
$ws.Range("E9").Value = 0.02929205403412615
$ws.Range("B10").Value = 30.21000000000001
$ws.Range("C10").Value = 30.26265525817871
$ws.Range("D10").Value = 0.05265525817870298
$ws.Range("E10").Value = 0.002772576213865867
$ws.Range("C11").Value = 30.34427452087402
$ws.Range("D11").Value = 0.1242745208740246
$ws.Range("E11").Value = 0.01544415653846837
$ws.Range("C12").Value = 30.36523246765137
$ws.Range("D12").Value = -0.01476753234862827
$ws.Range("E12").Value = 0.0002180800116677822
$ws.Range("C13").Value = 30.50981712341309
$ws.Range("D13").Value = 0.06981712341308821
$ws.Range("E13").Value = 0.00487443072167839
$ws.Range("C14").Value = 30.43168449401855
$ws.Range("D14").Value = -0.04831550598144929
$ws.Range("E14").Value = 0.002334388118243462
$ws.Range("C15").Value = 30.43945121765137
$ws.Range("D15").Value = -0.2505487823486305
$ws.Range("E15").Value = 0.06277469233638144
$ws.Range("C16").Value = 30.4976634979248
$ws.Range("D16").Value = -0.2523365020751953
$ws.Range("E16").Value = 0.06367371027954505
$ws.Range("C17").Value = 30.65968132019043
$ws.Range("D17").Value = -0.280318679809568
$ws.Range("E17").Value = 0.07857856225017913
$ws.Range("C18").Value = 30.76798057556152
$ws.Range("D18").Value = -0.1820194244384794
$ws.Range("E18").Value = 0.03313107087291531
$ws.Range("C19").Value = 31.12988090515137
$ws.Range("D19").Value = 0.1098809051513712
$ws.Range("E19").Value = 0.01207381331688463
$ws.Range("C20").Value = 31.24157524108887
$ws.Range("D20").Value = 0.1215752410888626
$ws.Range("E20").Value = 0.01478053924581508
$ws.Range("C21").Value = 31.3477611541748
$ws.Range("D21").Value = 0.06776115417480355
$ws.Range("E21").Value = 0.004591574015101497
$ws.Range("C22").Value = 31.23878288269043
$ws.Range("D22").Value = -0.1412171173095658
$ws.Range("E22").Value = 0.01994227422122366
$ws.Range("C23").Value = 31.39718437194824
$ws.Range("D23").Value = -0.1828156280517561
$ws.Range("E23").Value = 0.03342155385995804
$ws.Range("B24").Value = 31.65000000000001
$ws.Range("C24").Value = 31.92793846130371
$ws.Range("D24").Value = 0.2779384613037053
$ws.Range("E24").Value = 0.07724978827187126
$ws.Range("C25").Value = 32.41188430786133
$ws.Range("D25").Value = 0.5318843078613327
$ws.Range("E25").Value = 0.2829009169491289
$ws.Range("C26").Value = 32.36374282836914
$ws.Range("D26").Value = 0.08374282836913949
$ws.Range("E26").Value = 0.007012861303263154
$ws.Range("C27").Value = 32.45510482788086
$ws.Range("D27").Value = 0.005104827880856533
$ws.Range("E27").Value = 0.0000260592676931702
$ws.Range("B28").Value = 32.84999999999999
$ws.Range("C28").Value = 32.71161270141602
$ws.Range("D28").Value = -0.1383872985839787
$ws.Range("E28").Value = 0.01915104440937127
$ws.Range("B29").Value = 32.90000000000001
$ws.Range("C29").Value = 32.94222640991211
$ws.Range("D29").Value = 0.04222640991210369
$ws.Range("E29").Value = 0.001783069694065009
$ws.Range("B30").Value = 33.09999999999999
$ws.Range("C30").Value = 32.91791915893555
$ws.Range("D30").Value = -0.1820808410644474
$ws.Range("E30").Value = 0.03315343268273657
$ws.Range("B31").Value = 33.40000000000001
$ws.Range("C31").Value = 33.66357040405273
$ws.Range("D31").Value = 0.2635704040527287
$ws.Range("E31").Value = 0.06946935789251867
$ws.Range("C32").Value = 33.69541549682617
$ws.Range("D32").Value = -0.004584503173830967
$ws.Range("E32").Value = 0.00002101766935086621
$ws.Range("B33").Value = 34.09999999999999
$ws.Range("C33").Value = 33.89573287963867
$ws.Range("D33").Value = -0.2042671203613224
$ws.Range("E33").Value = 0.04172505646070699
$ws.Range("B34").Value = 34.40000000000001
$ws.Range("C34").Value = 34.42705917358398
$ws.Range("D34").Value = 0.02705917358397869
$ws.Range("E34").Value = 0.0007321988750478901
$ws.Range("B35").Value = 34.90000000000001
$ws.Range("C35").Value = 35.06875228881836
$ws.Range("D35").Value = 0.1687522888183537
$ws.Range("E35").Value = 0.02847733498143306
$ws.Range("C36").Value = 35.66774368286133
$ws.Range("D36").Value = 0.367743682861331
$ws.Range("E36").Value = 0.1352354162844152
$ws.Range("C37").Value = 35.96606826782227
$ws.Range("D37").Value = 0.2660682678222628
$ws.Range("E37").Value = 0.07079232314193935
$ws.Range("C38").Value = 35.86572647094727
$ws.Range("D38").Value = -0.4342735290527315
$ws.Range("E38").Value = 0.1885934980359137
$ws.Range("C39").Value = 36.4784049987793
$ws.Range("D39").Value = -0.3215950012207003
$ws.Range("E39").Value = 0.1034233448101422
$ws.Range("C40").Value = 37.2089958190918
$ws.Range("D40").Value = -0.09100418090820028
$ws.Range("E40").Value = 0.008281760942772444
$ws.Range("B41").Value = 37.90000000000001
$ws.Range("C41").Value = 37.99245071411133
$ws.Range("D41").Value = 0.09245071411132244
$ws.Range("E41").Value = 0.008547134539693474
$ws.Range("C42").Value = 38.43264389038086
$ws.Range("D42").Value = -0.06735610961914062
$ws.Range("E42").Value = 0.004536845503025688
$ws.Range("B43").Value = 38.90000000000001
$ws.Range("C43").Value = 39.07284927368164
$ws.Range("D43").Value = 0.1728492736816349
$ws.Range("E43").Value = 0.02987687141226874
$ws.Range("B44").Value = 39.40000000000001
$ws.Range("C44").Value = 39.60612869262695
$ws.Range("D44").Value = 0.2061286926269474
$ws.Range("E44").Value = 0.04248903792409458
$ws.Range("B45").Value = 39.90000000000001
$ws.Range("C45").Value = 39.64629745483398
$ws.Range("D45").Value = -0.2537025451660213
$ws.Range("E45").Value = 0.06436498142371708
$ws.Range("B46").Value = 40.09999999999999
$ws.Range("C46").Value = 39.90364456176758
$ws.Range("D46").Value = -0.1963554382324162
$ws.Range("E46").Value = 0.03855545812344421
$ws.Range("B47").Value = 40.59999999999999
$ws.Range("C47").Value = 40.41103744506836
$ws.Range("D47").Value = -0.1889625549316349
$ws.Range("E47").Value = 0.03570684716629115
$ws.Range("B48").Value = 40.90000000000001
$ws.Range("C48").Value = 40.64894485473633
$ws.Range("D48").Value = -0.2510551452636776
$ws.Range("E48").Value = 0.06302868596336623
$ws.Range("B49").Value = 41.20000000000001
$ws.Range("C49").Value = 41.17435073852539
$ws.Range("D49").Value = -0.02564926147461932
$ws.Range("E49").Value = 0.000657884614193391
$ws.Range("C50").Value = 41.48037338256836
$ws.Range("D50").Value = -0.01962661743164062
$ws.Range("E50").Value = 0.0003852041118079796
$ws.Range("C51").Value = 42.27613067626953
$ws.Range("D51").Value = 0.4761306762695341
$ws.Range("E51").Value = 0.2267004208848839
$ws.Range("C52").Value = 0.0128058624267382
$ws.Range("E52").Value = 1.993847182585832
$ws.Range("E53").Value = 0.03987694365171665
